$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")
$ws.Activate()

# Row 22: new "PC" request
$ws.Range("A22").Value = "Client"
$ws.Range("B22").Value = "Serveur"
$ws.Range("C22").Value = "PC"
$ws.Range("D22").Value = "idClasseChoisie"
$ws.Range("E22").Value = "Le client vérouille son choix de classe."

# Row 23: new "PO" response
$ws.Range("A23").Value = "Serveur"
$ws.Range("B23").Value = "Client"
$ws.Range("C23").Value = "PO"
$ws.Range("E23").Value = "Indique au client que son choix de classe est vérrouillé (aucune erreur n'est survenue lors de la génération automatique du placement de départ)"

$ws.Rows.Item(23).RowHeight = 60

# Update selection to reflect the scrolled-down view state
$ws.Range("A23").Select()

